{"js": "// Lattice-multiplication worksheet: regenerate all 15 problems in the\n// 5x3 table with a new set of multiplicand/multiplier pairs, per the\n// commit's refreshed \"output generated at c986bee\" data.\n//\n// Each cell holds a single run of 5 lines (joined with manual line\n// breaks): \"AB x CD\", the multiplier's digits spaced out, a divider,\n// and the two rows of the multiplicand's digits used to seed the\n// lattice grid. We rebuild each cell's text from the (A, B) pair so\n// the four derived lines always stay internally consistent, and we\n// re-insert each cell as OOXML so the \"  d    d\" / \"  ----\" lines keep\n// the xml:space=\"preserve\" the original markup uses for its\n// space-padded lines (plain w:t is fine for the other lines, which\n// never start/end with whitespace).\n\nconst pairs = [\n  [26, 29], [36, 58], [94, 63],\n  [50, 86], [99, 49], [97, 50],\n  [15, 23], [27, 29], [52, 54],\n  [26, 70], [49, 56], [99, 62],\n  [51, 41], [29, 12], [75, 13],\n];\n\nfunction cellOoxml(a, b) {\n  const aStr = String(a).padStart(2, \"0\");\n  const bStr = String(b).padStart(2, \"0\");\n  const line0 = a + \" x \" + b;\n  const line1 = \"  \" + bStr[0] + \"    \" + bStr[1];\n  const line2 = \"  ----\";\n  const line3 = aStr[0] + \"|    |\";\n  const line4 = aStr[1] + \"|    |\";\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:xml=\"http://www.w3.org/XML/1998/namespace\">' +\n    \"<w:body><w:p><w:r>\" +\n    '<w:rPr><w:sz w:val=\"32\"/></w:rPr>' +\n    \"<w:t>\" + line0 + \"</w:t><w:br/>\" +\n    '<w:t xml:space=\"preserve\">' + line1 + \"</w:t><w:br/>\" +\n    '<w:t xml:space=\"preserve\">' + line2 + \"</w:t><w:br/>\" +\n    \"<w:t>\" + line3 + \"</w:t><w:br/>\" +\n    \"<w:t>\" + line4 + \"</w:t>\" +\n    \"</w:r></w:p></w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nfor (const row of table.rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nlet idx = 0;\nfor (const row of table.rows.items) {\n  for (const cell of row.cells.items) {\n    const [a, b] = pairs[idx++];\n    cell.body.clear();\n    cell.body.insertOoxml(cellOoxml(a, b), Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Lattice-multiplication worksheet: regenerate all 15 problems in the\n# 5x3 table with a new set of multiplicand/multiplier pairs, per the\n# commit's refreshed \"output generated at c986bee\" data.\n#\n# Each cell holds a single run of 5 lines (separated by manual line\n# breaks, char 11 / vertical-tab, the same marker Range.Text uses for\n# <w:br/>): \"AB x CD\", the multiplier's digits spaced out, a divider,\n# and the two rows of the multiplicand's digits used to seed the\n# lattice grid. We rebuild each cell's text from the (A, C) pair so\n# the four derived lines always stay internally consistent.\n#\n# NB: string pieces are combined with the \"-f\" format operator rather\n# than \"+\" \u2014 PowerShell's \"+\" silently does numeric addition when both\n# operands look like numbers (e.g. \"  2    \" + \"9\" -> 11), which would\n# corrupt the padded digit rows.\n\n$pairs = @(\n  @(26, 29), @(36, 58), @(94, 63),\n  @(50, 86), @(99, 49), @(97, 50),\n  @(15, 23), @(27, 29), @(52, 54),\n  @(26, 70), @(49, 56), @(99, 62),\n  @(51, 41), @(29, 12), @(75, 13)\n)\n\n$VT = [char]11\n\nfunction Get-CellText([int]$a, [int]$b) {\n    $aStr = \"{0:D2}\" -f $a\n    $bStr = \"{0:D2}\" -f $b\n    $line0 = \"{0} x {1}\" -f $a, $b\n    $line1 = \"  {0}    {1}\" -f $bStr.Substring(0,1), $bStr.Substring(1,1)\n    $line2 = \"  ----\"\n    $line3 = \"{0}|    |\" -f $aStr.Substring(0,1)\n    $line4 = \"{0}|    |\" -f $aStr.Substring(1,1)\n    return (\"{0}{5}{1}{5}{2}{5}{3}{5}{4}\" -f $line0, $line1, $line2, $line3, $line4, $VT)\n}\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $pair = $pairs[$idx]\n        $idx++\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = Get-CellText $pair[0] $pair[1]\n    }\n}\n"}
